{"js": "// Update the date label and the 25 \"two-digit \u00d7 two-digit\" answer cells\n// to the new day's generated values. Each old value is unique in the\n// document, so a targeted search + Replace keeps existing run formatting\n// (font/size) intact.\nconst pairs = [\n  [\"2025-02-14 Friday\", \"2025-02-15 Saturday\"],\n  [\"13\u00d763=819\", \"39\u00d719=741\"],\n  [\"48\u00d779=3792\", \"16\u00d777=1232\"],\n  [\"80\u00d741=3280\", \"53\u00d776=4028\"],\n  [\"59\u00d715=885\", \"34\u00d712=408\"],\n  [\"66\u00d738=2508\", \"55\u00d716=880\"],\n  [\"22\u00d734=748\", \"53\u00d772=3816\"],\n  [\"80\u00d731=2480\", \"59\u00d782=4838\"],\n  [\"41\u00d777=3157\", \"78\u00d776=5928\"],\n  [\"33\u00d743=1419\", \"85\u00d777=6545\"],\n  [\"61\u00d748=2928\", \"72\u00d732=2304\"],\n  [\"95\u00d735=3325\", \"27\u00d777=2079\"],\n  [\"19\u00d717=323\", \"73\u00d789=6497\"],\n  [\"86\u00d797=8342\", \"35\u00d752=1820\"],\n  [\"70\u00d718=1260\", \"68\u00d775=5100\"],\n  [\"43\u00d759=2537\", \"55\u00d724=1320\"],\n  [\"44\u00d772=3168\", \"88\u00d793=8184\"],\n  [\"71\u00d791=6461\", \"21\u00d741=861\"],\n  [\"25\u00d762=1550\", \"53\u00d755=2915\"],\n  [\"33\u00d763=2079\", \"24\u00d784=2016\"],\n  [\"93\u00d730=2790\", \"51\u00d784=4284\"],\n  [\"67\u00d729=1943\", \"35\u00d747=1645\"],\n  [\"64\u00d798=6272\", \"67\u00d773=4891\"],\n  [\"13\u00d765=845\", \"97\u00d737=3589\"],\n  [\"25\u00d739=975\", \"45\u00d747=2115\"],\n  [\"27\u00d765=1755\", \"15\u00d764=960\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and the 25 \"two-digit \u00d7 two-digit\" answer cells\n# to the new day's generated values. Each old value occurs exactly once\n# in the document, so Find/Replace (wdReplaceAll) on each pair leaves\n# the surrounding run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-02-14 Friday\", \"2025-02-15 Saturday\"),\n    @(\"13\u00d763=819\", \"39\u00d719=741\"),\n    @(\"48\u00d779=3792\", \"16\u00d777=1232\"),\n    @(\"80\u00d741=3280\", \"53\u00d776=4028\"),\n    @(\"59\u00d715=885\", \"34\u00d712=408\"),\n    @(\"66\u00d738=2508\", \"55\u00d716=880\"),\n    @(\"22\u00d734=748\", \"53\u00d772=3816\"),\n    @(\"80\u00d731=2480\", \"59\u00d782=4838\"),\n    @(\"41\u00d777=3157\", \"78\u00d776=5928\"),\n    @(\"33\u00d743=1419\", \"85\u00d777=6545\"),\n    @(\"61\u00d748=2928\", \"72\u00d732=2304\"),\n    @(\"95\u00d735=3325\", \"27\u00d777=2079\"),\n    @(\"19\u00d717=323\", \"73\u00d789=6497\"),\n    @(\"86\u00d797=8342\", \"35\u00d752=1820\"),\n    @(\"70\u00d718=1260\", \"68\u00d775=5100\"),\n    @(\"43\u00d759=2537\", \"55\u00d724=1320\"),\n    @(\"44\u00d772=3168\", \"88\u00d793=8184\"),\n    @(\"71\u00d791=6461\", \"21\u00d741=861\"),\n    @(\"25\u00d762=1550\", \"53\u00d755=2915\"),\n    @(\"33\u00d763=2079\", \"24\u00d784=2016\"),\n    @(\"93\u00d730=2790\", \"51\u00d784=4284\"),\n    @(\"67\u00d729=1943\", \"35\u00d747=1645\"),\n    @(\"64\u00d798=6272\", \"67\u00d773=4891\"),\n    @(\"13\u00d765=845\", \"97\u00d737=3589\"),\n    @(\"25\u00d739=975\", \"45\u00d747=2115\"),\n    @(\"27\u00d765=1755\", \"15\u00d764=960\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
